# Trade #16 closed at 2026-02-17 15:18:00 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Summary sheet - update aggregate stats
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.78   # Current Capital
$summary.Range("B4").Value = -0.22     # Total P&L $
$summary.Range("B6").Value = 16        # Total Trades
$summary.Range("B8").Value = 7         # Losing Trades
$summary.Range("B9").Value = 25        # Win Rate %

# ---------------------------------------------------------------------------
# 2. Strategy Status sheet - update MarketMaking row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.78
$status.Range("D4").Value = 16
$status.Range("E4").Value = -0.22
$status.Range("F4").Value = -0.22
$status.Range("G4").Value = 25

# ---------------------------------------------------------------------------
# 3. Append the new trade row (#16) to both "All Trades" and "MarketMaking"
# ---------------------------------------------------------------------------
$tradeSheets = @("All Trades", "MarketMaking")

foreach ($sheetName in $tradeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A17").Value = 16

    # Dates/times are stored as plain text in this workbook (see existing
    # rows), so force a text number format before assignment to stop the
    # COM layer from coercing the string into a date/time serial value.
    $ws.Range("B17").NumberFormat = "@"
    $ws.Range("B17").Value = "2026-02-17"

    $ws.Range("C17").NumberFormat = "@"
    $ws.Range("C17").Value = "15:17:53"

    $ws.Range("D17").Value = "MarketMaking"
    $ws.Range("E17").Value = "DOWN"
    $ws.Range("F17").Value = 0.11
    $ws.Range("G17").Value = 0.09
    $ws.Range("H17").Value = "CLOSED"
    $ws.Range("I17").Value = -18.1818
    $ws.Range("J17").Value = -0.02
    $ws.Range("K17").Value = 99.78
    $ws.Range("L17").Value = 0
    $ws.Range("M17").Value = 0
    $ws.Range("N17").Value = 0.6
    $ws.Range("O17").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P17").Value = "early_exit"
    $ws.Range("Q17").Value = 0.18
}
